$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.000238706212329176,0.000179029659246882,0.000537088977740646,0.00173062003938653,0.996121024049651,0,0.998150026854449,0.000775795190069822,0.000059676553082294,0.138091543832428,0.00125320761472817,0.99504684609417,0.131944858864952,0.998866145491436,0.998746792385272,0.000537088977740646,0.311153547771081,0.00190964969863341,0.166139523781106,0.000417735871576058,0.998746792385272,0.0949453959539297,0.00238706212329176),
    @(0.981082532672913,0.999283881363012,0.999045175150683,0.0659425911559348,0.000238706212329176,0.000358059318493764,0.000954824849316704,0.000119353106164588,0.000179029659246882,0.000775795190069822,0.99844840961986,0.000179029659246882,0.000119353106164588,0.000417735871576058,0.000417735871576058,0.916631855344035,0.631079548845259,0.00745956913528675,0.315330906486841,0.999164528256848,0.00029838276541147,0.00029838276541147,0.000059676553082294),
    @(0.000417735871576058,0.000238706212329176,0.000238706212329176,0.0019693262517157,0.00310318076027929,0.000119353106164588,0.000417735871576058,0.998687115832189,0.000059676553082294,0.855344035328519,0.000119353106164588,0.00435638837500746,0.864176165184699,0.00059676553082294,0.000656442083905234,0.00029838276541147,0.0327027510890971,0.00179029659246882,0.167810467267411,0.000238706212329176,0.000656442083905234,0.902011099838873,0.996896819239721),
    @(0.0164110520976308,0.000238706212329176,0.000119353106164588,0.928925225278988,0.00029838276541147,0.999522587575342,0.000358059318493764,0.000358059318493764,0.999701617234589,0.000417735871576058,0.000119353106164588,0.000119353106164588,0.00029838276541147,0.000059676553082294,0.000119353106164588,0.0786536969624635,0.0211851763442144,0.988780808020529,0.342304708480038,0.000179029659246882,0.000238706212329176,0.000358059318493764,0.000537088977740646)
)

for ($i = 0; $i -lt 4; $i++) {
    for ($j = 0; $j -lt 23; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value = $data[$i][$j]
    }
}
